$wb = $excel.ActiveWorkbook

# "汽車" (Car) is the 3rd worksheet in this workbook (xl/worksheets/sheet3.xml).
$ws = $wb.Worksheets.Item("汽車")

# Make room for 7 new trailing columns (H:N) that carry the same
# bookkeeping fields used on the other property sheets (property_category,
# category, date, legislator_name, legislator_id, source_file, index).
$ws.Range("H1:N1").EntireColumn.Insert()

# --- Row 1: header labels -------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: data (B2:G2 unchanged, H2:N2 newly populated) -----------------
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Force text formatting first so the yyyy-mm-dd-looking string is kept
# literal instead of being auto-converted into a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-27"
$ws.Range("K2").Value = "李貴敏"
$ws.Range("L2").Value = 1739
$ws.Range("M2").Value = "tmp59331"
$ws.Range("N2").Value = 31
